$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 3625
$ws.Range("C3").Value = 3625
$ws.Range("C4").Value = 3938
$ws.Range("C5").Value = 3938
$ws.Range("C6").Value = 3938
$ws.Range("C7").Value = 3938
$ws.Range("C8").Value = 3938
$ws.Range("C9").Value = 3938
$ws.Range("C10").Value = 3938
$ws.Range("C11").Value = 3938
$ws.Range("C12").Value = 3938
